# Apply minor updates to the 焦煤288口岸监管区总库存_月度数据 monthly data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 216.2

$ws.Range("B3").Value = 218.2
$ws.Range("C3").Value = 246

$ws.Range("C4").Value = 336.5
$ws.Range("C5").Value = 351
$ws.Range("C6").Value = 391.6
$ws.Range("C7").Value = 394
$ws.Range("C8").Value = 370.8
$ws.Range("C9").Value = 397.3
$ws.Range("C10").Value = 439.7
$ws.Range("C11").Value = 436
$ws.Range("C12").Value = 363.6
$ws.Range("C13").Value = 413
$ws.Range("C14").Value = 384.6

$ws.Range("C16").Value = 357.8
$ws.Range("C17").Value = 358.5
$ws.Range("C18").Value = 304.2
$ws.Range("C19").Value = 324.9
$ws.Range("C20").Value = 307.8
$ws.Range("C21").Value = 338.5
$ws.Range("C22").Value = 274.9
$ws.Range("C23").Value = 272.5
$ws.Range("C24").Value = 216.7
